# "Generate Report for Archive"
#
# Localization status moved on from handoff: every cell that used to read
# "Ready for handoff" now reads "In Translation" (Overview sheet's per-locale
# status columns, plus each locale sheet's own Status column). Shrinking that
# text lets the Status-ish columns narrow, so their custom widths come down
# to match.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# New narrower custom width (XML `width` attribute target ~13.4101845877511).
# Excel quantizes ColumnWidth to pixel granularity and stores it as
# ColumnWidth + 5/6, so back the COM value out from the desired stored width.
$targetStoredWidth = 13.4101845877511
$newColumnWidth = $targetStoredWidth - (5.0 / 6.0)

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($addr in @("E2", "F2", "E3", "F3", "E4", "F4")) {
    $wsOverview.Range($addr).Value = $newStatus
}
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# --- Locale sheets: zh-cn and de-de, Status column C ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in @("C2", "C3", "C4")) {
        $ws.Range($addr).Value = $newStatus
    }
    $ws.Columns.Item(3).ColumnWidth = $newColumnWidth
}
